# Swap the "detalhar diária" (TC2) and "cancelar diária" (TC3) step/result
# text between the two test case blocks, per the commit's reordering of
# shared strings (v1.2.1 -> v1.2.3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Current (before) values:
#   B18 = "Beneficiário Clica em detalhar diária."
#   D18 = "SYSTEM Apresenta a tela de Detalhar Diárias"
#   B25 = "Beneficiário Clica em cancelar diária."
#   D25 = "SYSTEM Apresenta a tela de Cancelar Solicitação de Diária"
#
# Desired (after) values: swap the detalhar/cancelar content between the
# TC2 row (18) and TC3 row (25), while keeping TC2/TC3 labels in place.

$ws.Range("B18").Value = "Beneficiário Clica em cancelar diária."
$ws.Range("D18").Value = "SYSTEM Apresenta a tela de Cancelar Solicitação de Diária"

$ws.Range("B25").Value = "Beneficiário Clica em detalhar diária."
$ws.Range("D25").Value = "SYSTEM Apresenta a tela de Detalhar Diárias"
